$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each crypto row.
# D values must remain TEXT (the source data uses "." both as decimal
# and thousands separator, e.g. "26.070.05"), so force text format,
# write the value, then restore the default "Normal" style so the
# cell keeps its original (unstyled) appearance.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.070.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.718.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3695"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.57%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.71%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3334"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07482"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9994"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.923"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.715.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001077"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06647"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9985"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.073"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26.005.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.478"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.487"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.328"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.894.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.973"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08521"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("E35").Value = "  +2.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.371"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.279"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06210"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02289"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2135"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.516"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6174"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9991"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.839"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5889"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.012"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07256"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.59%  "

